# Update the calendar data "Readings"/"Reading Equivalences" header cells so
# the underlying export format used by convert_xlsx_to_yaml_calendar matches
# (spaces replaced with underscores, e.g. "Reading_Notes").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Reading_Notes"
$ws.Range("D1").Value = "Reading_Equivalences"
